$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$rows = @(
    @("2026-02-01", "18:19:11", "18:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom"),
    @("2026-02-01", "18:19:23", "18:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom"),
    @("2026-02-01", "18:19:30", "18:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom"),
    @("2026-02-01", "18:19:35", "18:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom"),
    @("2026-02-01", "18:19:43", "18:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom")
)

$startRow = 29
$endRow = $startRow + $rows.Count - 1

# Column A holds date-like text (e.g. "2026-02-01"); Excel would otherwise
# auto-convert it to a date serial. Force text entry, then clear the
# resulting number-format override so the cell keeps the default style.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$ws.Range("A$startRow`:A$endRow").ClearFormats()
